$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two new columns P and Q,
# copying the existing bold/centered/bordered header style from O1.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Add the new P and Q data columns (rows 2-25), all filled with 2.
$ws.Range("P2:Q25").Value = 2

# Swap values in columns I, K and M for data rows 2-25:
# I: 1 -> 2, K: 2 -> 1, M: 1 -> 2
$ws.Range("I2:I25").Value = 2
$ws.Range("K2:K25").Value = 1
$ws.Range("M2:M25").Value = 2
